$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (values that Excel will not misinterpret as numbers)
$plainUpdates = @{
    'D2' = '23.495.56'
    'E2' = '  -0.81%  '
    'D3' = '1.648.73'
    'E3' = '  -0.12%  '
    'E4' = '  +0.29%  '
    'E5' = '  +0.30%  '
    'E6' = '  -1.51%  '
    'E7' = '  -0.15%  '
    'E10' = '  -1.33%  '
    'E11' = '  -2.38%  '
    'E12' = '  +0.30%  '
    'E13' = '  -2.33%  '
    'E14' = '  -2.00%  '
    'E15' = '  -0.19%  '
    'E16' = '  -2.67%  '
    'D17' = '1.657.95'
    'E17' = '  +0.45%  '
    'E18' = '  +0.51%  '
    'E19' = '  -0.14%  '
    'E20' = '  -0.03%  '
    'E21' = '  -2.04%  '
    'E22' = '  +0.21%  '
    'E23' = '  -1.40%  '
    'D24' = '23.522.76'
    'E24' = '  -0.65%  '
    'E25' = '  -0.68%  '
    'E26' = '  -6.39%  '
    'E27' = '  -1.87%  '
    'E28' = '  -0.13%  '
    'E29' = '  +0.28%  '
    'E30' = '  -1.40%  '
    'D31' = '1.831.62'
    'E31' = '  -0.07%  '
    'E32' = '  +0.90%  '
    'B33' = 'FraxShare'
    'C33' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E33' = '  +4.14%  '
    'B34' = 'WEMIXTOKEN'
    'C34' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'E34' = '  +3.03%  '
    'E35' = '  -6.94%  '
    'E36' = '  -2.95%  '
    'E37' = '  -0.99%  '
    'B38' = 'InternetComputer(DFINITY)'
    'C38' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E38' = '  -2.10%  '
    'B39' = 'Algorand'
    'C39' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E39' = '  -3.39%  '
    'E40' = '  +2.55%  '
    'E41' = '  -3.28%  '
    'E42' = '  -2.40%  '
    'E43' = '  -1.47%  '
    'E44' = '  -2.06%  '
    'B45' = 'Decentraland'
    'C45' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E45' = '  -1.45%  '
    'B46' = 'Frax'
    'C46' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'E46' = '  +0.25%  '
    'E47' = '  -3.30%  '
    'E48' = '  -1.47%  '
    'E49' = '  -3.18%  '
    'E50' = '  -0.31%  '
    'E51' = '  -3.25%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Price updates whose text would otherwise be auto-coerced to a number by Excel
# (e.g. "1.001" -> 1.001, losing the trailing zero). Force the cell to Text first
# so the literal string is preserved exactly as in the source data.
$textForceUpdates = @{
    'D4' = '1.001'
    'D5' = '1.001'
    'D6' = '299.74'
    'D7' = '0.3800'
    'D8' = '0.3563'
    'D9' = '50.21'
    'D10' = '0.08095'
    'D11' = '1.220'
    'D12' = '1.001'
    'D13' = '22.03'
    'D14' = '6.402'
    'D15' = '7.391'
    'D16' = '0.00001198'
    'D18' = '97.27'
    'D19' = '0.06962'
    'D20' = '6.750'
    'D21' = '17.29'
    'D22' = '1.000'
    'D23' = '12.41'
    'D25' = '2.510'
    'D26' = '2.901'
    'D27' = '20.92'
    'D28' = '152.67'
    'D29' = '5.213'
    'D30' = '132.89'
    'D32' = '6.925'
    'D33' = '11.99'
    'D34' = '2.126'
    'D36' = '0.02729'
    'D37' = '0.08740'
    'D38' = '5.954'
    'D39' = '0.2433'
    'D40' = '13.15'
    'D41' = '0.06816'
    'D42' = '0.6904'
    'D43' = '1.316'
    'D44' = '15.54'
    'D45' = '0.6415'
    'D46' = '1.001'
    'D47' = '2.263'
    'D48' = '3.922'
    'D49' = '0.07737'
    'D50' = '127.75'
    'D51' = '1.155'
}
foreach ($addr in $textForceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForceUpdates[$addr]
}
